$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

for ($row = 2; $row -le 14; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value2 = $cell.Value2 + 1
}
